$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.38 = 29977.86 pesos`n✅ 29977.86 pesos = 7.35 = 970.28 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- tasas: update the usdt/bs/pesos rate table ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 135.5
$ws2.Range("O10").Value = 4062
$ws2.Range("N12").Value = 4079.5
$ws2.Range("O12").Value = 132.04
